{"js": "// The \"Seasonal Planning with Weather Data\" bullet had its body split across\n// three runs with stray w:proofErr (gramStart/gramEnd) markers wrapping the\n// word \"sowing\" in the middle. Collapse it back into a single clean run\n// (no proofErr markers) with the full, unchanged sentence text.\nconst fullText =\n  \": A farmer planning for the next season inputs upcoming weather \" +\n  \"forecasts and soil data into the app. The app recommends alternative \" +\n  \"crops or sowing dates based on projected rainfall, helping them adapt \" +\n  \"to changing weather patterns and maintain productivity.\";\n\nconst sentenceMatches = context.document.body.search(fullText, { matchCase: true });\nsentenceMatches.load(\"items\");\nawait context.sync();\n\nif (sentenceMatches.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the seasonal-planning sentence, found \" +\n      sentenceMatches.items.length\n  );\n}\n\nconst sentenceRange = sentenceMatches.items[0];\n\n// Build a minimal OOXML package with one run carrying the merged text, in\n// the same font as the surrounding content, and replace the matched range\n// with it. Unlike insertText(...,\"Replace\"), this removes the now-orphaned\n// w:proofErr elements left behind between the old runs.\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p><w:r><w:rPr>\" +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  \"</w:rPr><w:t>\" +\n  fullText +\n  \"</w:t></w:r></w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nsentenceRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Right below that bullet sits a stray paragraph whose only content is the\n// leftover text \"4o\". Remove that run's text, leaving the empty paragraph.\nconst strayMatches = context.document.body.search(\"4o\", { matchCase: true });\nstrayMatches.load(\"items\");\nawait context.sync();\n\nif (strayMatches.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the stray '4o' text, found \" +\n      strayMatches.items.length\n  );\n}\n\nstrayMatches.items[0].delete();\nawait context.sync();\n", "ps1": "# The \"Seasonal Planning with Weather Data\" bullet text was split across\n# several runs, with the word \"sowing\" wrapped in a pair of grammar-check\n# w:proofErr markers. Collapse the sentence back into plain, unmarked text by\n# finding the full sentence and replacing it with itself: Find/Replace across\n# a multi-run match rewrites the range as a single clean run and drops the\n# now-orphaned proofErr markers.\n$d = $word.ActiveDocument\n\n$sentence = \": A farmer planning for the next season inputs upcoming weather forecasts and soil data into the app. The app recommends alternative crops or sowing dates based on projected rainfall, helping them adapt to changing weather patterns and maintain productivity.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$found = $find.Execute(\n    $sentence,  # FindText\n    $false,     # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $sentence,  # ReplaceWith\n    2           # Replace (wdReplaceAll)\n)\nif (-not $found) {\n    throw \"Seasonal-planning sentence not found\"\n}\n\n# Just below that bullet sits a stray paragraph whose only content is the\n# leftover text \"4o\". Clear that run's text, leaving the empty paragraph.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$found2 = $find2.Execute(\n    \"4o\",   # FindText\n    $true,  # MatchCase\n    $false, # MatchWholeWord\n    $false, # MatchWildcards\n    $false, # MatchSoundsLike\n    $false, # MatchAllWordForms\n    $true,  # Forward\n    1,      # Wrap (wdFindContinue)\n    $false, # Format\n    \"\",     # ReplaceWith\n    2       # Replace (wdReplaceAll)\n)\nif (-not $found2) {\n    throw \"Stray '4o' text not found\"\n}\n"}
